$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new text labels (set in this order so shared strings line up) ---
$ws.Range("A23").Value = "Dev Main + Database + docker containers"
$ws.Range("D24").Value = "load1"
$ws.Range("A24").Value = "NGINX load balancer for k3s servers"
$ws.Range("B19").Value = "Don’t use this worker:"

# --- Row 17: bump E17 / H17 by 1, highlight E17 in yellow ---
$ws.Range("E17").Value = 7
$ws.Range("E17").Interior.Color = 65535
$ws.Range("H17").Value = 5

# --- Row 18: bump E18 / H18 by 1, highlight E18 in yellow ---
$ws.Range("E18").Value = 7
$ws.Range("E18").Interior.Color = 65535
$ws.Range("H18").Value = 5

# --- Row 19: mark this worker as unused, zero out its resources ---
$ws.Range("E19").ClearContents()
$ws.Range("E19").Interior.Color = 65535
$ws.Range("H19").Value = 0
$ws.Range("H19").Interior.Color = 65535
$ws.Range("M19").Value = 0
$ws.Range("M19").Interior.Color = 65535

# --- Row 23: rework numbers for the Dev Main box ---
$ws.Range("E23").Interior.Color = 65535
$ws.Range("H23").Interior.Color = 65535
$ws.Range("K23").Value = 0
$ws.Range("K23").Interior.Color = 65535
$ws.Range("M23").Value = 200
$ws.Range("M23").Interior.Color = 65535

# --- Row 24 (new): add the NGINX load balancer / load1 entry ---
$ws.Range("P20").Copy() | Out-Null
$ws.Range("A24").PasteSpecial(-4122) | Out-Null
$ws.Range("E24").Value = 4
$ws.Range("E24").Interior.Color = 65535
$ws.Range("F24").Formula = "=F23-E24"
$ws.Range("G24").Value = "dev1"
$ws.Range("H24").Value = 2
$ws.Range("H24").Interior.Color = 65535
$ws.Range("I24").Formula = "=I23-H24"
$ws.Range("K24").Value = 25
$ws.Range("K24").Interior.Color = 65535
$ws.Range("L24").Formula = "=L23-K24"

# --- cosmetic: last selected cell ---
$ws.Range("E19").Select() | Out-Null
